$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update floodmedia column (H) for all data rows from "NA" to "None"
$ws.Range("H2:H27").Value = "None"

# Row heights: rows 2 and 4-27 change from 16 to 15 (row 3 was already 15)
foreach ($r in 2,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27) {
    $ws.Rows($r).RowHeight = 15
}

# Update the selection to H2:H27, active cell H2
$ws.Range("H2:H27").Select()
